# Update NATMI Mfng-Notch1 LR-pair sheet with recomputed values from new TPM input.
# (commit: "update scripts wuth new tpm")
#
# The underlying per-cell TPM table that feeds this LR-pairs export changed,
# which changes ligand/receptor expression aggregates for the "ECs" and
# "MuSCs" clusters (the "FAPs" cluster aggregates are unaffected), and every
# downstream specificity / edge-weight column that derives from them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> { column letter -> new value }
$updates = @{
    2  = @{ G = 8.893517000000001;  H = 26.680551;  I = 0.9082944842335181; J = 0.9082944842335181;
            M = 38.55267666666666;  N = 115.65803;  O = 0.5758151725879548; P = 0.5758151725879548;
            Q = 342.8688853305034;  R = 3085.81996797453; S = 0.5230097451996106; T = 0.5230097451996106 }
    3  = @{ G = 8.893517000000001;  H = 26.680551;  I = 0.9082944842335181; J = 0.9082944842335181;
            O = 0.08021535714867321; P = 0.08021535714867323;
            Q = 47.76419830749201;  R = 429.8777847674281; S = 0.07285916644896159; T = 0.07285916644896159 }
    4  = @{ G = 8.893517000000001;  H = 26.680551;  I = 0.9082944842335181; J = 0.9082944842335181;
            M = 23.02986166666667;  N = 69.089585;  O = 0.3439694702633719; P = 0.3439694702633719;
            Q = 204.8164662401483; R = 1843.348196161335; S = 0.3124255725849458; T = 0.3124255725849459 }
    5  = @{ I = 0.04237443292342908; J = 0.04237443292342909;
            M = 38.55267666666666;  N = 115.65803;  O = 0.5758151725879548; P = 0.5758151725879548;
            Q = 15.99577541773666;  R = 143.96197875963; S = 0.02439984140712103; T = 0.02439984140712103 }
    6  = @{ I = 0.04237443292342908; J = 0.04237443292342909;
            O = 0.08021535714867321; P = 0.08021535714867323;
            S = 0.00339908027092536; T = 0.003399080270925362 }
    7  = @{ I = 0.04237443292342908; J = 0.04237443292342909;
            M = 23.02986166666667;  N = 69.089585;  O = 0.3439694702633719; P = 0.3439694702633719;
            Q = 9.555250814531666;  R = 85.997257330785; S = 0.01457551124538269; T = 0.01457551124538269 }
    8  = @{ E = 3; F = 1; G = 0.4830226666666667; H = 1.449068; I = 0.04933108284305281; J = 0.04933108284305281;
            M = 38.55267666666666;  N = 115.65803;  O = 0.5758151725879548; P = 0.5758151725879548;
            Q = 18.62181669067111; R = 167.59635021604; S = 0.02840558598122315; T = 0.02840558598122315 }
    9  = @{ E = 3; F = 1; G = 0.4830226666666667; H = 1.449068; I = 0.04933108284305281; J = 0.04933108284305281;
            O = 0.08021535714867321; P = 0.08021535714867323;
            Q = 2.594158243322667; R = 23.347424189904; S = 0.003957110428786267; T = 0.003957110428786268 }
    10 = @{ E = 3; F = 1; G = 0.4830226666666667; H = 1.449068; I = 0.04933108284305281; J = 0.04933108284305281;
            M = 23.02986166666667;  N = 69.089585;  O = 0.3439694702633719; P = 0.3439694702633719;
            Q = 11.12394519519778; R = 100.11550675678; S = 0.01696838643304339; T = 0.01696838643304339 }
}

foreach ($rowNum in $updates.Keys) {
    $cols = $updates[$rowNum]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$rowNum").Value = $cols[$col]
    }
}
